$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("P4").Value = 20215070055
$ws.Range("P5").Value = "Muhammed Ali Harmancı"
$ws.Range("P6").Value = "Yönetim Bilişim Sistemleri"
$ws.Range("G6").Formula = "=SUM(D4:D17)"
$ws.Range("D18").Value = 1484
Write-Host "Done"
